# Updates the cryptos price/volume table (Sheet1, rows 2-51) to the refreshed
# values from the scheduled GitHub Actions data pull.
# Column D ("Price") values that look purely numeric are written with a leading
# apostrophe so Excel stores them as text (quotePrefix) instead of silently
# converting them to numbers - this matches how the source data/sheet stores them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.742.86'
$ws.Range("E2").Value = '  +0.49%  '

$ws.Range("D3").Value = '''1.594.65'
$ws.Range("E3").Value = '  -0.19%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '''210.08'
$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("D6").Value = '''0.501'
$ws.Range("E6").Value = '  -0.19%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '''22.41'
$ws.Range("E8").Value = '  -0.13%  '

$ws.Range("E9").Value = '  +0.20%  '

$ws.Range("D10").Value = '''0.0592'
$ws.Range("E10").Value = '  +0.25%  '

$ws.Range("E11").Value = '  -0.69%  '

$ws.Range("D12").Value = '''1.821.46'
$ws.Range("E12").Value = '  -0.30%  '

$ws.Range("D13").Value = '''1.590.36'
$ws.Range("E13").Value = '  -0.38%  '

$ws.Range("D14").Value = '''3.84'
$ws.Range("E14").Value = '  -0.52%  '

$ws.Range("E15").Value = '  -1.59%  '

$ws.Range("D16").Value = '''27.748.29'
$ws.Range("E16").Value = '  +0.47%  '

$ws.Range("D17").Value = '''63.44'
$ws.Range("E17").Value = '  -0.56%  '

$ws.Range("D18").Value = '''219.01'
$ws.Range("E18").Value = '  +0.21%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0697'
$ws.Range("E19").Value = '  +0.05%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '''7.36'
$ws.Range("E20").Value = '  -1.16%  '

$ws.Range("E21").Value = '  +0.02%  '

$ws.Range("E22").Value = '  -1.01%  '

$ws.Range("D23").Value = '''9.79'
$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("E24").Value = '  -1.99%  '

$ws.Range("D25").Value = '''153.77'
$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("D26").Value = '''7.15'
$ws.Range("E26").Value = '  +6.26%  '

$ws.Range("E27").Value = '  +0.08%  '

$ws.Range("E28").Value = '  +0.81%  '

$ws.Range("E29").Value = '  -0.37%  '

$ws.Range("E30").Value = '  +0.57%  '

$ws.Range("D31").Value = '''0.0474'
$ws.Range("E31").Value = '  +1.33%  '

$ws.Range("E32").Value = '  -1.90%  '

$ws.Range("D33").Value = '''1.382.28'
$ws.Range("E33").Value = '  +0.56%  '

$ws.Range("E34").Value = '  +0.58%  '

$ws.Range("E35").Value = '  -0.76%  '

$ws.Range("D36").Value = '''0.971'
$ws.Range("E36").Value = '  +1.06%  '

$ws.Range("E37").Value = '  +0.84%  '

$ws.Range("D38").Value = '''0.0170'
$ws.Range("E38").Value = '  +2.74%  '

$ws.Range("E39").Value = '  +0.18%  '

$ws.Range("D40").Value = '''0.829'
$ws.Range("E40").Value = '  +1.05%  '

$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D42").Value = '''0.984'
$ws.Range("E42").Value = '  -0.53%  '

$ws.Range("D43").Value = '''64.60'
$ws.Range("E43").Value = '  +0.88%  '

$ws.Range("D44").Value = '''2.17'
$ws.Range("E44").Value = '  +3.46%  '

$ws.Range("D45").Value = '''1.76'
$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("D46").Value = '''5.26'
$ws.Range("E46").Value = '  -0.50%  '

$ws.Range("D47").Value = '''1.731.99'
$ws.Range("E47").Value = '  -0.35%  '

$ws.Range("D48").Value = '''85.78'
$ws.Range("E48").Value = '  -2.46%  '

$ws.Range("E49").Value = '  +5.16%  '

$ws.Range("D50").Value = '''0.0967'
$ws.Range("E50").Value = '  -0.54%  '

$ws.Range("E51").Value = '  -0.27%  '
